$wb = $excel.ActiveWorkbook

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 12:16:22"
$wsOverview.Range("G5").Value = "2016-08-24 12:16:22"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-24 12:16:18"
$wsZh.Range("K3").Value = "2016-08-24 12:16:35"
$wsZh.Range("E5").Value = "mt"
$wsZh.Range("H5").Value = "2016-08-24 12:16:18"
$wsZh.Range("K5").Value = "2016-08-24 12:16:35"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-24 12:16:22"
$wsDe.Range("K3").Value = "2016-08-24 12:16:43"
$wsDe.Range("E5").Value = "mt"
$wsDe.Range("H5").Value = "2016-08-24 12:16:22"
$wsDe.Range("K5").Value = "2016-08-24 12:16:43"
